$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (pushes existing data rows down by one)
$ws.Rows.Item(2).Insert()

# Populate the new row with the latest shipment data.
# Columns C, E and F look like plain numbers (leading zeros / long digit
# strings), so force those specific cells to text first - otherwise Excel's
# smart-typing would silently convert them to numeric values and drop the
# leading zeros. A / B / D are already alphanumeric and don't need this.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("A2").Value = "TLLU4140075"
$ws.Range("B2").Value = "TRITON"
$ws.Range("C2").Value = "00016"
$ws.Range("D2").Value = "DJLAXA3787663"
$ws.Range("E2").Value = "7075348470"
$ws.Range("F2").Value = "082900024555"

# Remove the two stale rows (previously rows 9 and 10, now shifted to 10 and 11)
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(10).Delete()
